# intraday koers sprint 038
# Insert two new diary entries (each preceded by a blank paragraph, and
# followed by a trailing blank paragraph) right after the paragraph that
# ends with "...bij koersen van indexen is deze link anders."

$d = $word.ActiveDocument

$apos = [char]0x2019

$anchorRange = $d.Content
$anchorRange.Find.Execute("16:33 Opgelost. Eindelijk. Bleek aan de paginalink te liggen, bij koersen van indexen is deze link anders.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorPara = $anchorRange.Paragraphs(1)

# 1) blank paragraph
$anchorPara.Range.InsertParagraphAfter()
$blank1 = $anchorPara.Next()

# 2) "Gaat lekker ... degiro ligt er werkelijk uit." paragraph
$blank1.Range.InsertParagraphAfter()
$entry1 = $blank1.Next()
$entry1.Range.InsertAfter("Gaat")
$entry1.Range.InsertAfter(" lekker. Degiro" + $apos + "s website lag er al eerder uit maar toen kon ik nog beleggen. Nu is de hele website en ook de beleggerssite onbereikbaar.")
$entry1.Range.InsertAfter(" Storing opgezocht, degiro ligt er werkelijk uit.")

# 3) blank paragraph
$entry1.Range.InsertParagraphAfter()
$blank2 = $entry1.Next()

# 4) "Ik zie Put RD maart 16 ..." paragraph
$blank2.Range.InsertParagraphAfter()
$entry2 = $blank2.Next()
$entry2.Range.InsertAfter("Ik zie Put RD maart 16 op 1,12 staan en het aandeel op 15,72. Hoef dus niks te doen. Bijna als resultaat 0 voor deze optie-exercitie.")

# 5) trailing blank paragraph
$entry2.Range.InsertParagraphAfter()
